$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60000
$ws.Range("C9").Value = "charge by rounding up to the nearest dollar"
$ws.Range("C3").Value = "security, food, greeting desk"
$ws.Range("C4").Value = "I will be maintenance"
$ws.Range("B9").Value = 0

$ws.Range("C4").Select()
